# Apply the crypto-price refresh update (GitHub Actions daily data pull).
# Only columns B (Coin), C (Link), D (Price) and E (Volume 1h) on rows 2-51
# are touched; rows 18/19 also swap their Coin/Link/Price/Volume content
# (ShibaInu now ranks above Dai).
#
# The "Price" column stores its values as plain text in the workbook (not
# numbers), even though most of them look numeric. Writing a numeric-looking
# string straight into .Value lets Excel auto-coerce it to a real number, so
# for those cells we enter the value the way a user forcing text entry would
# (a leading apostrophe) and then reset the cell style to "Normal" so no
# stray Text number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  D = '30.671.75';     E = '  -0.48%  ' },
    @{ Row = 3;  D = '1.923.27';      E = '  -0.69%  ' },
    @{ Row = 4;  D = '0.9987';        E = '  -0.19%  ' },
    @{ Row = 5;  D = '241.48';        E = '  -0.73%  ' },
    @{ Row = 6;  D = '0.9986';        E = '  -0.14%  ' },
    @{ Row = 7;  D = '0.4786';        E = '  -1.92%  ' },
    @{ Row = 8;  D = '0.2885';        E = '  -2.03%  ' },
    @{ Row = 9;  D = '0.06780';       E = '  -1.48%  ' },
    @{ Row = 10; D = '19.66';         E = '  +1.95%  ' },
    @{ Row = 11; D = '104.54';        E = '  -0.12%  ' },
    @{ Row = 12;                     E = '  -0.03%  ' },
    @{ Row = 13; D = '1.934.11';      E = '  -0.05%  ' },
    @{ Row = 14; D = '5.265';         E = '  -1.32%  ' },
    @{ Row = 15; D = '0.6805';        E = '  -2.77%  ' },
    @{ Row = 16; D = '295.82';        E = '  +8.41%  ' },
    @{ Row = 17; D = '30.687.67';     E = '  -0.38%  ' },
    @{ Row = 18; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000007575'; E = '  -1.89%  ' },
    @{ Row = 19; B = 'Dai';      C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai';       D = '0.9998';       E = '  -0.02%  ' },
    @{ Row = 20; D = '12.84';         E = '  -1.87%  ' },
    @{ Row = 21; D = '5.492';         E = '  -2.37%  ' },
    @{ Row = 22; D = '0.9984';        E = '  -0.29%  ' },
    @{ Row = 23; D = '6.386';         E = '  -2.10%  ' },
    @{ Row = 24; D = '9.495';         E = '  -3.26%  ' },
    @{ Row = 25; D = '168.22';        E = '  +1.80%  ' },
    @{ Row = 26; D = '19.77';         E = '  +0.70%  ' },
    @{ Row = 27; D = '2.115' },
    @{ Row = 28; D = '1.395';         E = '  +0.30%  ' },
    @{ Row = 29; D = '0.1004';        E = '  -3.19%  ' },
    @{ Row = 30; D = '4.611' },
    @{ Row = 31; D = '1.522';         E = '  -2.17%  ' },
    @{ Row = 32; D = '4.312';         E = '  -1.55%  ' },
    @{ Row = 33; D = '0.04784';       E = '  -2.06%  ' },
    @{ Row = 34; D = '0.7350';        E = '  -3.00%  ' },
    @{ Row = 35; D = '1.121';         E = '  -2.45%  ' },
    @{ Row = 36; D = '2.710';         E = '  -0.77%  ' },
    @{ Row = 37; D = '0.01928';       E = '  -3.91%  ' },
    @{ Row = 38; D = '2.628';         E = '  -1.11%  ' },
    @{ Row = 39; D = '6.411';         E = '  -1.25%  ' },
    @{ Row = 40; D = '75.27';         E = '  -5.97%  ' },
    @{ Row = 41; D = '1.994';         E = '  -4.27%  ' },
    @{ Row = 42; D = '0.8644';        E = '  -4.36%  ' },
    @{ Row = 43; D = '106.41';        E = '  -1.54%  ' },
    @{ Row = 44; D = '0.4307';        E = '  -2.77%  ' },
    @{ Row = 45; D = '0.9986';        E = '  -0.18%  ' },
    @{ Row = 46; D = '7.505';         E = '  -3.54%  ' },
    @{ Row = 47; D = '981.51';        E = '  -2.28%  ' },
    @{ Row = 48; D = '0.1213';        E = '  -2.62%  ' },
    @{ Row = 49; D = '34.96';         E = '  -3.13%  ' },
    @{ Row = 50; D = '8.921';         E = '  -3.28%  ' },
    @{ Row = 51; D = '0.05815';       E = '  +0.65%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { Set-TextValue $ws.Cells.Item($r, 4) $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}
